$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "B2" = 13.4015549233451
    "C2" = 6.11502654280375
    "E2" = 11.02526790736374
    "F2" = 16.86991607391233
    "G2" = 38.8828758685366
    "H2" = 16.8871261757979
    "I2" = 26.20550565969513
    "K2" = 10.88976906574689
    "L2" = 9.970671698646735
    "N2" = 20.14734012728668
    "B3" = 13.18567586768959
    "C3" = 6.041231258576162
    "E3" = 11.02905030871177
    "F3" = 15.89584955866808
    "G3" = 38.97400841493282
    "H3" = 16.94038787306201
    "I3" = 26.29789182702483
    "K3" = 10.74732747256997
    "L3" = 9.959120296097726
    "N3" = 20.21319274475371
    "B4" = 13.05421562322823
    "C4" = 5.994678695329648
    "E4" = 11.03321825506279
    "F4" = 15.26997757108491
    "G4" = 39.04097569701341
    "H4" = 16.97580167411344
    "I4" = 26.35942141664308
    "K4" = 10.66114578053583
    "L4" = 9.95379551555083
    "N4" = 20.25549693736497
    "B5" = 13.00098942836069
    "C5" = 5.975403827462253
    "E5" = 11.03538152614433
    "F5" = 15.00819731993403
    "G5" = 39.07102369296896
    "H5" = 16.99091442369053
    "I5" = 26.3857018323407
    "K5" = 10.62638893579598
    "L5" = 9.952072213318669
    "N5" = 20.27320800011668
    "B6" = 12.99217415936688
    "C6" = 5.972185107227752
    "E6" = 11.03576882644433
    "F6" = 14.96433081551593
    "G6" = 39.0761794046566
    "H6" = 16.99346502545126
    "I6" = 26.3901385098393
    "K6" = 10.62064070857563
    "L6" = 9.951813087768782
    "N6" = 20.27617744379844
    "B7" = 13.0534963060932
    "C7" = 5.9944199692057
    "E7" = 11.03324554689365
    "F7" = 15.26647399323137
    "G7" = 39.04136978138612
    "H7" = 16.9760027316916
    "I7" = 26.35977095991334
    "K7" = 10.66067551378595
    "L7" = 9.953770463762272
    "N7" = 20.25573388292451
    "B8" = 13.32693590694099
    "C8" = 6.089844416390566
    "E8" = 11.02618951681757
    "F8" = 16.53996406344768
    "G8" = 38.91200678429124
    "H8" = 16.9049279580457
    "I8" = 26.23636246548793
    "K8" = 10.84041501308356
    "L8" = 9.96632317514627
    "N8" = 20.16965876851581
    "B9" = 13.86871264355566
    "C9" = 6.266779678781358
    "E9" = 11.02695898836653
    "F9" = 19.00274580682531
    "G9" = 38.74613587819498
    "H9" = 16.78707554656292
    "I9" = 26.03254064089143
    "K9" = 11.20118143443381
    "L9" = 10.00486472763095
    "N9" = 20.01564029302036
    "B10" = 14.26607803768896
    "C10" = 6.390128721360526
    "E10" = 11.03637081601564
    "F10" = 20.67494806633232
    "G10" = 38.67834208772643
    "H10" = 16.71363233203923
    "I10" = 25.90615578054807
    "K10" = 11.4688240788691
    "L10" = 10.04151916000462
    "N10" = 19.91139529250402
    "B11" = 14.44587248967711
    "C11" = 6.444713939875925
    "E11" = 11.0425576965682
    "F11" = 21.3917225636224
    "G11" = 38.65933623471048
    "H11" = 16.68307948644117
    "I11" = 25.85375159873905
    "K11" = 11.59062598638763
    "L11" = 10.05996671368971
    "N11" = 19.86588656115339
    "B12" = 14.51375042182671
    "C12" = 6.465157391829121
    "E12" = 11.04517285016187
    "F12" = 21.65686569030329
    "G12" = 38.65384674153283
    "H12" = 16.67192114298235
    "I12" = 25.83464086328859
    "K12" = 11.63671414639799
    "L12" = 10.06720362921846
    "N12" = 19.8489271581068
    "B13" = 14.4991419079786
    "C13" = 6.460764726926056
    "E13" = 11.04459754729916
    "F13" = 21.60004134736742
    "G13" = 38.65495297427506
    "H13" = 16.67430598807794
    "I13" = 25.83872404679891
    "K13" = 11.62679048784384
    "L13" = 10.06563391940815
    "N13" = 19.85256751374565
    "B14" = 14.45146135688469
    "C14" = 6.446400420427797
    "E14" = 11.04276739944107
    "F14" = 21.4136618050453
    "G14" = 38.65885035813407
    "H14" = 16.6821532355115
    "I14" = 25.85216463446851
    "K14" = 11.59441864107547
    "L14" = 10.06055708664727
    "N14" = 19.86448582158265
    "B15" = 14.42222683216311
    "C15" = 6.437572124363617
    "E15" = 11.04168179184578
    "F15" = 21.29868154950795
    "G15" = 38.66146015466006
    "H15" = 16.68701348797284
    "I15" = 25.8604929711348
    "K15" = 11.57458406655913
    "L15" = 10.05747998081523
    "N15" = 19.87182174330561
    "B16" = 14.25430276041654
    "C16" = 6.386530100266046
    "E16" = 11.03600469089924
    "F16" = 20.62722412089977
    "G16" = 38.6798228255916
    "H16" = 16.71568656605718
    "I16" = 25.90968307429138
    "K16" = 11.46086136286742
    "L16" = 10.04034893589762
    "N16" = 19.9144077751491
    "B17" = 14.15099008034059
    "C17" = 6.354821171269869
    "E17" = 11.03300885338268
    "F17" = 20.20408069597325
    "G17" = 38.69412325525731
    "H17" = 16.73400866601941
    "I17" = 25.94116429738535
    "K17" = 11.39107820055615
    "L17" = 10.03029117779945
    "N17" = 19.94102192374925
    "B18" = 14.09147999151809
    "C18" = 6.336439730791044
    "E18" = 11.03146518311025
    "F18" = 19.95656407809801
    "G18" = 38.70346206590098
    "H18" = 16.74481592027132
    "I18" = 25.95975035054745
    "K18" = 11.35094796820195
    "L18" = 10.02467337263653
    "N18" = 19.95650980076327
    "B19" = 14.07131797729438
    "C19" = 6.33019170319066
    "E19" = 11.03097339284494
    "F19" = 19.87204792380568
    "G19" = 38.70681508584219
    "H19" = 16.74852123234343
    "I19" = 25.96612546262521
    "K19" = 11.33736311730725
    "L19" = 10.02280009904402
    "N19" = 19.96178470580442
    "B20" = 14.16199742751808
    "C20" = 6.35821152547742
    "E20" = 11.03330920394478
    "F20" = 20.24955283636154
    "G20" = 38.69248565191204
    "H20" = 16.73203041649566
    "I20" = 25.93776349133193
    "K20" = 11.39850630560912
    "L20" = 10.03134457021606
    "N20" = 19.93817016960327
    "B21" = 14.46547241021326
    "C21" = 6.450625776455253
    "E21" = 11.04329758189122
    "F21" = 21.46857628470577
    "G21" = 38.65765921677614
    "H21" = 16.67983714030681
    "I21" = 25.84819688427418
    "K21" = 11.60392833803995
    "L21" = 10.06204148860879
    "N21" = 19.8609777060352
    "B22" = 14.66257916851855
    "C22" = 6.509698808118721
    "E22" = 11.0514119286108
    "F22" = 22.22866616901552
    "G22" = 38.64485327999403
    "H22" = 16.64812376497127
    "I22" = 25.79393668607801
    "K22" = 11.73795712101265
    "L22" = 10.08356610899955
    "N22" = 19.81212309071585
    "B23" = 14.55751368354417
    "C23" = 6.478293969581383
    "E23" = 11.04693659638738
    "F23" = 21.82633154458858
    "G23" = 38.65077551059991
    "H23" = 16.66483019877992
    "I23" = 25.82250445774563
    "K23" = 11.66645798189715
    "L23" = 10.07194551688718
    "N23" = 19.83805219084493
    "B24" = 14.157021356587
    "C24" = 6.356679217883464
    "E24" = 11.03317285878212
    "F24" = 20.22900810905287
    "G24" = 38.69322253247798
    "H24" = 16.73292393076753
    "I24" = 25.93929947863553
    "K24" = 11.39514809017497
    "L24" = 10.03086781871994
    "N24" = 19.93945886488543
    "B25" = 13.72197795521119
    "C25" = 6.220050662217811
    "E25" = 11.02519241755807
    "F25" = 18.34778573295695
    "G25" = 38.7815478384278
    "H25" = 16.81665120232427
    "I25" = 26.08358327154827
    "K25" = 11.10295227773517
    "L25" = 9.992961891050346
    "N25" = 20.05573447111166
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}

Write-Output "Updated $($values.Count) cells"
